$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 607 (shifts the existing rows 607:648 down to 608:649,
# growing the used range from D648 to D649).
$ws.Rows.Item(607).Insert()

# The date column (A) stores plain text like "2026/12/29", not real dates.
# Assigning a date-shaped string directly would make Excel auto-coerce it
# into a date serial, so force Text formatting first, then strip the
# resulting style override back off so the cell matches its neighbours
# (which carry no explicit style at all).
$ws.Range("A607").NumberFormat = "@"
$ws.Range("A607").Value = "2026/01/12"
$ws.Range("A607").ClearFormats()

$ws.Range("B607").Value = "月"
$ws.Range("C607").Value = 7
$ws.Range("D607").Value = 18
